# "Added some corrections in my report"
#
# 1) The deck's date placeholder (master + every custom layout) was
#    refreshed from 18/05/2015 to 24/05/2015.
# 2) On the single slide, the "+" / "-" correction labels on the two
#    outer axis arrows were swapped (the outer-right label was wrongly
#    marked "+" and the outer-left one wrongly marked "-").

$p = $ppt.ActivePresentation

# --- 1) Refresh the "date updated automatically" placeholder everywhere ---
function Update-DateText($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = "24/05/2015"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateText $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateText $layout.Shapes
}

# --- 2) Swap the "+"/"-" correction signs on slide 1 ---
$slide = $p.Slides.Item(1)

$plusBox = $slide.Shapes.Item("25 CuadroTexto")
$minusBox = $slide.Shapes.Item("27 CuadroTexto")

$plusBox.TextFrame.TextRange.Text = "-"
$minusBox.TextFrame.TextRange.Text = "+"
